$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the column headers in row 1:
#      "<name>_old" -> "<name>_FV2410"
#      "<name>_new" -> "<name>_FV2504"
#    (the "diff" header in column K is left untouched)
# ---------------------------------------------------------------------------
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $base = $val.Substring(0, $val.Length - 4)
            $cell.Value = $base + "_FV2410"
        } elseif ($val.EndsWith("_new")) {
            $base = $val.Substring(0, $val.Length - 4)
            $cell.Value = $base + "_FV2504"
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into a native Excel Table ("Table1") so the
#    (now renamed) headers also drive the table's column names.
#    1 = xlSrcRange, $false = no LinkSource, 1 = xlYes (has headers)
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U80")
$tbl = $ws.ListObjects.Add(1, $tableRange, $false, 1)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3. Freeze the header row so it stays visible while scrolling.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
